# CS554 EuroTeam Proposal - add a new slide explaining the use case diagram.
#
# The new slide is inserted right before the existing "Use case" slide
# (which was slide 9), pushing it and everything after it down by one
# position. It uses the same custom layout ("Diapo classique" /
# slideLayout3.xml) as the neighbouring slides, with a centred title
# placeholder and a body text placeholder.

$p = $ppt.ActivePresentation

# Slide 9 ("Use case") is on the "Diapo classique" layout - reuse that
# same CustomLayout object for the freshly inserted slide.
$refSlide = $p.Slides.Item(9)
$layout = $refSlide.CustomLayout

$newSlide = $p.Slides.AddSlide(9, $layout)

# --- Title placeholder -------------------------------------------------
$title = $newSlide.Shapes.Item(1)
$title.Name = "Title 1"
$title.TextFrame.TextRange.Text = "Use case diagram explained"

# --- Body placeholder ----------------------------------------------------
$body = $newSlide.Shapes.Item(2)
$body.Name = "Text Placeholder 2"
$body.TextFrame.TextRange.Text = "Most of the time FDRI works automatically. However, spaceship crew and flight control can manually control the system.`rThese interactions are represented in the following use case diagram.`r"
